$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("C1").Value = 4
$ws.Range("D1").Value = 9
$ws.Range("A2").Value = 8
$ws.Range("B2").Value = 5
$ws.Range("A3").Value = 7
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 2
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = 8
$ws.Range("A6").Value = 10
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 5

# Update selection to C7
$ws.Range("C7").Select()

# Update window view yWindow 540 -> 560
$excel.ActiveWindow.Top = 560
